$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook used to track a generic "Property" concept per-sheet; this
# tab is being folded into the unified DataNode/DataTable/Entity model, so
# rename it accordingly.
$ws.Name = "DataNode"

# Leave the cursor parked on the cell the editor ended the session on.
$ws.Range("H33").Select() | Out-Null
